$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.531.56"
$ws.Range("E2").Value = "  -2.84%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.756.73"
$ws.Range("E3").Value = "  -2.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.98"
$ws.Range("E5").Value = "  -0.55%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.40%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4444"
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3698"
$ws.Range("E8").Value = "  -1.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.93"
$ws.Range("E9").Value = "  +0.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07668"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("E11").Value = "  -3.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.62"
$ws.Range("E13").Value = "  -4.19%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.440"
$ws.Range("E14").Value = "  -2.95%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.153"
$ws.Range("E15").Value = "  -2.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.758.78"
$ws.Range("E16").Value = "  -2.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.34"
$ws.Range("E17").Value = "  +11.85%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001071"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06256"
$ws.Range("E19").Value = "  -7.98%  "

# Row 20
$ws.Range("E20").Value = "  +0.19%  "

# Row 21
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.168"
$ws.Range("E22").Value = "  -2.27%  "

# Row 23
$ws.Range("E23").Value = "  -3.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.574.99"
$ws.Range("E24").Value = "  -2.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.55"
$ws.Range("E25").Value = "  -2.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.312"
$ws.Range("E26").Value = "  -4.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.55"
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.13"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.286"
$ws.Range("E29").Value = "  -3.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.957.67"
$ws.Range("E30").Value = "  -2.47%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.81"
$ws.Range("E31").Value = "  -3.69%  "

# Row 32
$ws.Range("E32").Value = "  -5.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.709"
$ws.Range("E33").Value = "  -1.89%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09205"
$ws.Range("E34").Value = "  -1.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.634"
$ws.Range("E35").Value = "  -9.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.57"
$ws.Range("E36").Value = "  +3.94%  "

# Row 37
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2160"
$ws.Range("E38").Value = "  -5.44%  "

# Row 39
$ws.Range("E39").Value = "  -3.99%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.049"
$ws.Range("E40").Value = "  -2.08%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6433"
$ws.Range("E41").Value = "  -2.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.170"
$ws.Range("E42").Value = "  -3.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.946"
$ws.Range("E43").Value = "  -2.42%  "

# Row 44
$ws.Range("E44").Value = "  +0.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.391"
$ws.Range("E45").Value = "  -4.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.72"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5967"
$ws.Range("E47").Value = "  -1.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.724"
$ws.Range("E48").Value = "  -1.83%  "

# Row 49
$ws.Range("E49").Value = "  -2.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.996"
$ws.Range("E50").Value = "  -1.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06891"
$ws.Range("E51").Value = "  -2.88%  "
